$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.501.99'
$ws.Range('E2').Value = '  -0.72%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.626.32'
$ws.Range('E3').Value = '  -0.65%  '
$ws.Range('E4').Value = '  +0.33%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '212.97'
$ws.Range('E5').Value = '  -0.13%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.502'
$ws.Range('E6').Value = '  +1.40%  '
$ws.Range('E7').Value = '  +0.28%  '
$ws.Range('E8').Value = '  -0.29%  '
$ws.Range('E9').Value = '  -1.52%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '18.79'
$ws.Range('E10').Value = '  -1.33%  '
$ws.Range('E11').Value = '  +0.53%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.854.02'
$ws.Range('E12').Value = '  -0.62%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.627.42'
$ws.Range('E13').Value = '  -0.43%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.14'
$ws.Range('E14').Value = '  +1.26%  '
$ws.Range('E15').Value = '  -0.93%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '65.04'
$ws.Range('E16').Value = '  +2.99%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '26.539.78'
$ws.Range('E17').Value = '  -0.64%  '
$ws.Range('E18').Value = '  -0.36%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '214.40'
$ws.Range('E19').Value = '  +2.37%  '
$ws.Range('E20').Value = '  +0.28%  '
$ws.Range('E21').Value = '  -1.05%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.25'
$ws.Range('E22').Value = '  +1.15%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.30'
$ws.Range('E23').Value = '  -1.24%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.08'
$ws.Range('E24').Value = '  +9.30%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '147.81'
$ws.Range('E25').Value = '  +0.80%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.01'
$ws.Range('E26').Value = '  +0.29%  '
$ws.Range('E27').Value = '  -0.46%  '
$ws.Range('E28').Value = '  +1.87%  '
$ws.Range('E29').Value = '  +0.60%  '
$ws.Range('E30').Value = '  -2.12%  '
$ws.Range('E31').Value = '  -1.19%  '
$ws.Range('E32').Value = '  +3.12%  '
$ws.Range('E33').Value = '  -0.27%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.243.61'
$ws.Range('E34').Value = '  +6.28%  '
$ws.Range('E35').Value = '  -0.36%  '
$ws.Range('E36').Value = '  +0.40%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0175'
$ws.Range('E37').Value = '  +4.25%  '
$ws.Range('E38').Value = '  +0.28%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.508'
$ws.Range('E39').Value = '  +0.51%  '
$ws.Range('E40').Value = '  -1.97%  '
$ws.Range('E41').Value = '  -1.93%  '
$ws.Range('E42').Value = '  +0.64%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.34'
$ws.Range('E43').Value = '  -0.94%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.764.31'
$ws.Range('E44').Value = '  -0.99%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '93.21'
$ws.Range('E45').Value = '  +0.91%  '
$ws.Range('E46').Value = '  +1.63%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '54.85'
$ws.Range('E47').Value = '  +0.00%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0₆0104'
$ws.Range('E48').Value = '  -0.47%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0510'
$ws.Range('E49').Value = '  -0.62%  '
$ws.Range('E50').Value = '  -0.54%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.49'
$ws.Range('E51').Value = '  -0.94%  '
